$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4067
$ws.Range("J62").Value = 4067
$ws.Range("L62").Value = 4067
$ws.Range("N62").Value = -5315

$ws.Range("H65").Value = 4067
$ws.Range("J65").Value = 4067
$ws.Range("L65").Value = 20335
$ws.Range("N65").Value = -26575

$ws.Range("H100").Value = 1245.25
$ws.Range("I100").Value = 751.2
$ws.Range("K100").Value = 751.2
$ws.Range("M100").Value = -210.2

$ws.Range("H103").Value = 1495.8
$ws.Range("I103").Value = 1495.8
$ws.Range("K103").Value = 4487.4
$ws.Range("M103").Value = -3901.4

$ws.Range("H106").Value = 66670464
$ws.Range("I106").Value = 111112440
$ws.Range("K106").Value = 111112440
$ws.Range("M106").Value = -111111809

$ws.Range("H116").Value = 33590880
$ws.Range("I116").Value = 27656414
$ws.Range("J116").Value = 45459816
$ws.Range("K116").Value = 27656414
$ws.Range("L116").Value = 45459816
$ws.Range("M116").Value = -27652972
$ws.Range("N116").Value = -45466700

$ws.Range("H132").Value = 9947.044
$ws.Range("I132").Value = 1525.25
$ws.Range("K132").Value = 4575.75
$ws.Range("M132").Value = -2045.75

$ws.Range("H133").Value = 93351.17999999999
$ws.Range("J133").Value = 93351.17999999999
$ws.Range("L133").Value = 93351.17999999999
$ws.Range("N133").Value = -103471.18

$ws.Range("H137").Value = 6669629
$ws.Range("I137").Value = 1154.3103
$ws.Range("J137").Value = 15878474
$ws.Range("K137").Value = 3462.9309
$ws.Range("L137").Value = 47635422
$ws.Range("M137").Value = -912.9309000000003
$ws.Range("N137").Value = -47640522

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3450.3066
$ws.Range("I32").Value = 1890.7377
$ws.Range("K32").Value = 1890.7377
$ws.Range("M32").Value = -1603.7377

$ws.Range("H102").Value = 572244.0600000001
$ws.Range("I102").Value = 1523446.1
$ws.Range("K102").Value = 1523446.1
$ws.Range("M102").Value = -1521824.1

$ws.Range("H122").Value = 5986.85
$ws.Range("I122").Value = 3605.8462
$ws.Range("J122").Value = 10408.714
$ws.Range("K122").Value = 10817.5386
$ws.Range("L122").Value = 31226.142
$ws.Range("M122").Value = -8367.5386
$ws.Range("N122").Value = -36126.142

$ws.Range("H132").Value = 15493.347
$ws.Range("I132").Value = 17189.027
$ws.Range("K132").Value = 51567.08099999999
$ws.Range("M132").Value = -49037.08099999999

$ws.Range("H139").Value = 86439.39999999999
$ws.Range("J139").Value = 68843.336
$ws.Range("L139").Value = 68843.336
$ws.Range("N139").Value = -79123.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1706.8948
$ws.Range("I86").Value = 1762.25
$ws.Range("J86").Value = 1666.6364
$ws.Range("K86").Value = 1762.25
$ws.Range("L86").Value = 1666.6364
$ws.Range("M86").Value = -639.25
$ws.Range("N86").Value = -3912.6364

$ws.Range("H89").Value = 1706.8948
$ws.Range("I89").Value = 1762.25
$ws.Range("J89").Value = 1666.6364
$ws.Range("K89").Value = 8811.25
$ws.Range("L89").Value = 8333.182000000001
$ws.Range("M89").Value = -3195.25
$ws.Range("N89").Value = -19565.182

$ws.Range("H105").Value = 1937.5
$ws.Range("I105").Value = 1110.2858
$ws.Range("K105").Value = 1110.2858
$ws.Range("M105").Value = 636.7141999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 500
$ws.Range("J8").Value = 500
$ws.Range("L8").Value = 500
$ws.Range("N8").Value = -780

$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").ClearContents()
$ws.Range("N25").Value = 0

$ws.Range("H31").Value = 3159.96
$ws.Range("I31").Value = 1327.6428
$ws.Range("J31").Value = 5492
$ws.Range("K31").Value = 1327.6428
$ws.Range("L31").Value = 5492
$ws.Range("M31").Value = -1032.6428
$ws.Range("N31").Value = -6082

$ws.Range("H34").Value = 3159.96
$ws.Range("I34").Value = 1327.6428
$ws.Range("J34").Value = 5492
$ws.Range("K34").Value = 1327.6428
$ws.Range("L34").Value = 5492
$ws.Range("M34").Value = -1125.6428
$ws.Range("N34").Value = -5896

$ws.Range("H62").Value = 28098.25
$ws.Range("I62").Value = 3300.8
$ws.Range("K62").Value = 3300.8
$ws.Range("M62").Value = -2676.8

$ws.Range("H65").Value = 28098.25
$ws.Range("I65").Value = 3300.8
$ws.Range("K65").Value = 16504
$ws.Range("M65").Value = -13384

$ws.Range("H103").Value = 52994.6
$ws.Range("J103").Value = 65368.25
$ws.Range("L103").Value = 65368.25
$ws.Range("N103").Value = -67712.25

$ws.Range("H132").Value = 30306020
$ws.Range("I132").Value = 33336022
$ws.Range("J132").Value = 5998
$ws.Range("K132").Value = 100008066
$ws.Range("L132").Value = 17994
$ws.Range("M132").Value = -100005536
$ws.Range("N132").Value = -23054

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1999
$ws.Range("I129").Value = 1298.8
$ws.Range("K129").Value = 3896.4
$ws.Range("M129").Value = 1103.6

$ws.Range("H133").Value = 28196
$ws.Range("I133").Value = 11187.5
$ws.Range("K133").Value = 33562.5
$ws.Range("M133").Value = -28502.5

$ws.Range("H137").Value = 4083737.8
$ws.Range("J137").Value = 4593510.5
$ws.Range("L137").Value = 13780531.5
$ws.Range("N137").Value = -13790731.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I70").Value = 3973654.2
$ws.Range("K70").Value = 3973654.2
$ws.Range("M70").Value = -3973384.2

$ws.Range("I73").Value = 3973654.2
$ws.Range("K73").Value = 3973654.2
$ws.Range("M73").Value = -3972718.2

$ws.Range("H102").Value = 12926.077
$ws.Range("I102").Value = 14004.333
$ws.Range("K102").Value = 14004.333
$ws.Range("M102").Value = -12382.333

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").ClearContents()
$ws.Range("N106").Value = 0

$ws.Range("H126").Value = 4269.0713
$ws.Range("I126").Value = 2757.1428
$ws.Range("J126").Value = 5781
$ws.Range("K126").Value = 8271.428400000001
$ws.Range("L126").Value = 17343
$ws.Range("M126").Value = -5801.428400000001
$ws.Range("N126").Value = -22283

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 8001813.5
$ws.Range("I16").Value = 9092220
$ws.Range("K16").Value = 9092220
$ws.Range("M16").Value = -9092050

$ws.Range("H46").Value = 5915.514
$ws.Range("I46").Value = 3431.125
$ws.Range("J46").Value = 6651.6294
$ws.Range("K46").Value = 3431.125
$ws.Range("L46").Value = 6651.6294
$ws.Range("M46").Value = -3243.125
$ws.Range("N46").Value = -7027.6294

$ws.Range("H132").Value = 4045.8386
$ws.Range("I132").Value = 3507.2632
$ws.Range("J132").Value = 4898.5835
$ws.Range("K132").Value = 10521.7896
$ws.Range("L132").Value = 14695.7505
$ws.Range("M132").Value = -7991.7896
$ws.Range("N132").Value = -19755.7505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 1166.6666
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 1500
$ws.Range("K4").Value = 1000
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = -887
$ws.Range("N4").Value = -1726

$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws.Range("H64").Value = 99996
$ws.Range("J64").Value = 99996
$ws.Range("L64").Value = 99996
$ws.Range("N64").Value = -100492

$ws.Range("H67").Value = 99996
$ws.Range("J67").Value = 99996
$ws.Range("L67").Value = 99996
$ws.Range("N67").Value = -101712

$ws.Range("H100").Value = 1821151.2
$ws.Range("I100").Value = 4002160
$ws.Range("K100").Value = 8004320
$ws.Range("M100").Value = -8003779

$ws.Range("H107").Value = 2637.1538
$ws.Range("I107").Value = 2862.5
$ws.Range("J107").Value = 1397.75
$ws.Range("K107").Value = 8587.5
$ws.Range("L107").Value = 4193.25
$ws.Range("M107").Value = -6667.5
$ws.Range("N107").Value = -8033.25

$ws.Range("H132").Value = 13230283
$ws.Range("I132").Value = 2138771.5
$ws.Range("K132").Value = 6416314.5
$ws.Range("M132").Value = -6413784.5

$ws.Range("H133").Value = 56763.332
$ws.Range("J133").Value = 56763.332
$ws.Range("L133").Value = 56763.332
$ws.Range("N133").Value = -66883.33199999999
